$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D1").Value = "Admin Url: https://djangotask.herokuapp.com/admin/"
